# Fix field validation and UI guide - update test Excel file to match
# the standard template format.
#
# Sheet "Input" gets a full column restructure:
#   - New set of headers (발주일자, 납기일자, 거래처명, 거래처 이메일,
#     납품처명, 납품처 이메일, 프로젝트명, 대분류, 중분류, 소분류,
#     품목명, 규격, 수량, 단가, 총금액, 비고) replacing the old 17-column
#     layout with a 16-column layout.
#   - Header row loses its bold/bordered "title" style.
#   - Data rows are re-populated/re-ordered to match the new headers.
# Sheets "갑지" / "을지" just lose their stray empty trailing "비고" cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Input
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Input")

# Wipe the old layout (values, styles, and the now-stray column Q) completely
# so the sheet can be rebuilt from scratch with the new column order/width.
$ws.UsedRange.Clear()

$headers = @("발주일자", "납기일자", "거래처명", "거래처 이메일", "납품처명", "납품처 이메일", "프로젝트명", "대분류", "중분류", "소분류", "품목명", "규격", "수량", "단가", "총금액", "비고")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Helper: write a value that looks like a date ("2025-09-13") as literal
# text instead of letting Excel auto-convert it to a date serial number.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 1) "2025-09-13"
Set-TextValue $ws.Cells.Item(2, 2) "2025-10-10"
$ws.Cells.Item(2, 3).Value = "이노메탈"
$ws.Cells.Item(2, 4).Value = "이노메탈@example.com"
$ws.Cells.Item(2, 5).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(2, 6).Value = "delivery@example.com"
$ws.Cells.Item(2, 7).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(2, 8).Value = "5. 운반비"
$ws.Cells.Item(2, 9).Value = "일반자재"
$ws.Cells.Item(2, 10).Value = "기타"
$ws.Cells.Item(2, 11).Value = "4월 운반비"
$ws.Cells.Item(2, 12).Value = "KS규격-1"
$ws.Cells.Item(2, 13).Value = 1
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0

# Row 3
Set-TextValue $ws.Cells.Item(3, 1) "2025-08-28"
Set-TextValue $ws.Cells.Item(3, 2) "2025-09-24"
$ws.Cells.Item(3, 3).Value = "이노메탈"
$ws.Cells.Item(3, 4).Value = "이노메탈@example.com"
$ws.Cells.Item(3, 5).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(3, 6).Value = "delivery@example.com"
$ws.Cells.Item(3, 7).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(3, 8).Value = "1. 원자재비"
$ws.Cells.Item(3, 9).Value = "4) ALUM. 창호"
$ws.Cells.Item(3, 10).Value = "B. 도장"
$ws.Cells.Item(3, 11).Value = "품목명 없음"
$ws.Cells.Item(3, 12).Value = "KS규격-2"
$ws.Cells.Item(3, 13).Value = 1307
$ws.Cells.Item(3, 14).Value = 2600
$ws.Cells.Item(3, 15).Value = 3738020

# ---------------------------------------------------------------------
# Sheets 2 & 3: 갑지 / 을지 - drop the stray empty "비고" cells at I2/I3
# ---------------------------------------------------------------------
foreach ($name in @("갑지", "을지")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Cells.Item(2, 9).ClearContents()
    $sheet.Cells.Item(3, 9).ClearContents()
}
